# Fix sheet "基金受益憑證" (funds, sheet #6):
#   - row 1 currently duplicates row 2's data instead of real column headers
#   - rows 2-9 are missing the trailing metadata columns that every other
#     sheet has (property_category, category, date, legislator_name,
#     legislator_id, source_file, index)
#   - a "dealer" column (D) needs a proper header label
#   - H5 is stored as text "t65312" instead of the numeric total 65312

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# ---- Row 1: real headers (copy format from the bold/boxed header cells) ----
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "dealer"
$ws.Range("E1").Value = "quantity"
$ws.Range("F1").Value = "face_value"
$ws.Range("G1").Value = "currency"
$ws.Range("H1").Value = "total"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

# new header cells (I1:O1) get the same bold/boxed style as the rest of row 1
$ws.Range("B1").Copy()
$ws.Range("I1:O1").PasteSpecial(-4122)

# ---- Fix H5: numeric total, not text "t65312" ----
$ws.Range("H5").Value = 65312

# ---- Rows 2-9: append metadata columns I:O ----
for ($r = 2; $r -le 9; $r++) {
    $idx = $ws.Cells.Item($r, 1).Value

    $ws.Cells.Item($r, 9).Value  = "fund"
    $ws.Cells.Item($r, 10).Value = "normal"
    $ws.Cells.Item($r, 11).Value = "2012-04-27"
    $ws.Cells.Item($r, 12).Value = "陳雪生"
    $ws.Cells.Item($r, 13).Value = 1751
    $ws.Cells.Item($r, 14).Value = "tmp5a001"
    $ws.Cells.Item($r, 15).Value = $idx

    # match the bordered "data row" style used by the rest of the row
    $ws.Range("B" + $r).Copy()
    $ws.Range("I" + $r + ":O" + $r).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
